$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Add the new "USB Cabel" product as row 22 (same layout/pattern as the
# other product rows above it): Name | Price | Qty | Link
# ---------------------------------------------------------------------
$productUrl = "https://www.conrad.de/de/p/garmin-usb-kabel-usb-2-0-usb-a-stecker-usb-mini-a-stecker-1-00-m-schwarz-010-10723-01-373148.html"

# A22 - product name
$ws.Range("A22").Value2 = "USB Cabel"
$ws.Range("A22").HorizontalAlignment = -4108
$ws.Range("A22").VerticalAlignment = -4108
$ws.Range("A22").WrapText = $False

# B22 - price, Euro currency format
$ws.Range("B22").Value2 = 6.49
$ws.Range("B22").NumberFormat = "[$€-2]\ #,##0.00"
$ws.Range("B22").HorizontalAlignment = -4108
$ws.Range("B22").VerticalAlignment = -4108
$ws.Range("B22").WrapText = $True

# C22 - quantity
$ws.Range("C22").Value2 = 1
$ws.Range("C22").HorizontalAlignment = -4108
$ws.Range("C22").VerticalAlignment = -4108
$ws.Range("C22").WrapText = $False

# D22 - link to the product page (display text mirrors the other rows,
# which show the raw URL as the hyperlink's visible text)
$ws.Hyperlinks.Add($ws.Range("D22"), $productUrl, "", "", $productUrl) | Out-Null
$ws.Range("D22").HorizontalAlignment = -4108
$ws.Range("D22").VerticalAlignment = -4108
$ws.Range("D22").WrapText = $True

# Row height to fit the wrapped link text, like the other product rows
$ws.Rows.Item(22).RowHeight = 108

# ---------------------------------------------------------------------
# Extend the total formula to include the new row
# ---------------------------------------------------------------------
$ws.Range("B31").Formula = "=SUM(B2:B22)"

# ---------------------------------------------------------------------
# Update the active selection to reflect the new last row
# ---------------------------------------------------------------------
$ws.Range("B32").Select() | Out-Null
